$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 1.08
$ws.Range("L2").Value = 1.05
$ws.Range("T2").Value = 1.13
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.06
$ws.Range("X2").Value = "2026-01-15 08:00:00"
$ws.Range("Z2").Value = 1.05
